$d = $word.ActiveDocument

$pairs = @(
    @("392÷3=", "281÷2="),
    @("855÷7=", "583÷2="),
    @("974÷5=", "236÷2="),
    @("512÷8=", "312÷4="),
    @("166÷3=", "368÷5="),
    @("122÷4=", "273÷3="),
    @("672÷4=", "966÷9="),
    @("975÷7=", "742÷9="),
    @("694÷3=", "173÷3="),
    @("701÷7=", "256÷5="),
    @("372÷8=", "534÷4="),
    @("903÷7=", "227÷8="),
    @("564÷5=", "342÷8="),
    @("988÷8=", "412÷8="),
    @("250÷5=", "540÷9="),
    @("524÷3=", "976÷4="),
    @("362÷8=", "646÷2="),
    @("309÷2=", "960÷7="),
    @("723÷6=", "312÷4="),
    @("514÷8=", "915÷8="),
    @("163÷3=", "460÷6="),
    @("905÷6=", "802÷7="),
    @("696÷5=", "105÷9="),
    @("918÷8=", "363÷8="),
    @("447÷8=", "744÷5=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
